$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.653.46'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -3.08%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.097.94'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.72%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.007'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.32%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '342.52'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.006'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.36%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5125'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.45%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4397'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.55%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '53.22'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.37%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.09145'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.36%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.170'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.35%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '24.80'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.22%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.097.30'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.46%  '
$ws.Range('E14').Value = '  -1.04%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.183'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.66%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '99.78'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.62%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001145'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.99%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.007'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.50%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '21.07'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +8.48%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.06640'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.27%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.006'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.38%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.181'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.75%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '29.698.11'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.58'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.86%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.305'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.26%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.343.16'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.76%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.82'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.67%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '162.25'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.65%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.524'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.54%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '132.46'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.14%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.130'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -5.06%  '
$ws.Range('E32').Value = '  -3.03%  '
$ws.Range('E33').Value = '  -1.58%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.153'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.28%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.960'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.49%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.043'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.15%  '
$ws.Range('E37').Value = '  -1.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02567'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.87%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06661'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.63%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2236'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.20%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '12.36'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.36%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6846'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.47%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.287'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.22%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6663'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.29%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.22'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.292'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.74%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.608'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.03%  '
$ws.Range('E48').Value = '  -2.66%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '81.74'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.77%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00000000333'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -7.75%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.163'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.42%  '
